# Fix the Cypher WHERE-clause queries on the "startup" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# B2: StatQuery cell - fix the ethnicity filter syntax
$b2 = "MATCH (ct:clinical_trial)<--(a:arm)<--(c:case)`n    WHERE c.ethnicity = `"UNKNOWN`" `nWITH DISTINCT c, a, ct`nRETURN `n    COALESCE(c.case_id, '') AS ``Case ID``,`n    COALESCE(ct.clinical_trial_designation, '') AS ``Trial Code``,`n    COALESCE(a.arm_id, '') AS ``Arm``,`n    COALESCE(a.arm_drug, '') AS ``Arm Treatment``,`n    COALESCE(c.disease, '') AS ``Diagnosis``,`n    COALESCE(c.gender, '') AS ``Gender``,`n    COALESCE(c.race, '') AS ``Race``,`n    COALESCE(c.ethnicity, '') AS ``Ethnicity``"
$ws.Range("B2").Value = $b2

# C2: dbExcel query cell - fix the duplicate WHERE and the ethnicity filter syntax
$c2 = "MATCH (s:specimen)-->(c:case)-->(:arm)-->(ct:clinical_trial)`n    WHERE c.ethnicity = `"UNKNOWN`" `nOPTIONAL MATCH (f:file)-->(:sequencing_assay)-->(:nucleic_acid)-->(s)`nRETURN `n    COUNT(DISTINCT f) AS number_of_files,`n    COUNT(DISTINCT c.case_id) AS number_of_cases,`n    COUNT(DISTINCT ct.clinical_trial_designation) AS number_of_trials"
$ws.Range("C2").Value = $c2

# Update the active selection to B3 as recorded by the saved workbook view
$ws.Range("B3").Select()

$wb.Save()
